# "added pre-requisite for Account Lockout"
#
# LoginCredentials sheet: the old row 3 (testlock / locktest123 / testlock123 /
# "Account Lockout Test") is replaced by two rows: row 3 keeps the
# testlock/locktest123 valid-login pair (relabelled "Test valid password" in
# column C instead of the old testlock123/"Account Lockout Test" split across
# C3/D3) and a brand new row 4 is added with testlock/testlock123 labelled
# "Test invalid password" -- the pre-requisite data row that will be used to
# repeatedly fail a login and trigger the account lockout.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LoginCredentials")

# Drop the old D3 label ("Account Lockout Test") -- it's no longer used.
$ws1.Range("D3").Value = ""

# Row 3 now just reuses the "Test valid password" label in column C.
$ws1.Range("C3").Value = "Test valid password"

# New row 4: the account-lockout pre-requisite data.
$ws1.Range("A4").Value = "testlock"
$ws1.Range("B4").Value = "testlock123"
$ws1.Range("C4").Value = "Test invalid password"

# Match formatting: C3 should look like the plain (non-centered) label cells
# (same as C2), and A4/B4 should look like the centered credential cells
# (same as A3/B3), while C4 stays like the plain label cells.
$ws1.Range("C2").Copy()
$ws1.Range("C3").PasteSpecial(-4122)

$ws1.Range("A3:B3").Copy()
$ws1.Range("A4").PasteSpecial(-4122)

$ws1.Range("C2").Copy()
$ws1.Range("C4").PasteSpecial(-4122)

$ws1.Range("E7").Select()

# CreateEmployee sheet: re-touch the hyperlinked e-mail cells so each gets
# its own relationship entry instead of all three sharing one.
$ws4 = $wb.Worksheets.Item("CreateEmployee")
$ws4.Range("M4").Value = "martink@testmail.com"
$ws4.Range("M3").Value = "martink@testmail.com"
$ws4.Range("M2").Value = "martink@testmail.com"
